$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $style
}

$ws.Range("D2").Value = "28.865.01"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.878.53"
$ws.Range("E3").Value = "  -1.90%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "325.07"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -0.10%  "
Set-TextValue $ws.Range("D7") "0.4611"
$ws.Range("E7").Value = "  -1.15%  "
Set-TextValue $ws.Range("D8") "0.3872"
$ws.Range("E8").Value = "  -2.09%  "
Set-TextValue $ws.Range("D9") "0.07845"
$ws.Range("E9").Value = "  -2.47%  "
Set-TextValue $ws.Range("D10") "0.9846"
$ws.Range("E10").Value = "  -3.34%  "
Set-TextValue $ws.Range("D11") "21.74"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "1.870.70"
$ws.Range("E12").Value = "  -5.03%  "
Set-TextValue $ws.Range("D13") "6.991"
$ws.Range("E13").Value = "  -2.62%  "
Set-TextValue $ws.Range("D14") "5.657"
$ws.Range("E14").Value = "  -2.58%  "
Set-TextValue $ws.Range("D15") "0.06990"
$ws.Range("E15").Value = "  +0.10%  "
Set-TextValue $ws.Range("D16") "88.04"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("E17").Value = "  -0.02%  "
Set-TextValue $ws.Range("D18") "0.000009954"
$ws.Range("E18").Value = "  -2.39%  "
Set-TextValue $ws.Range("D19") "16.90"
$ws.Range("E19").Value = "  -2.76%  "
Set-TextValue $ws.Range("D20") "1.002"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "28.861.14"
$ws.Range("E21").Value = "  -1.34%  "
Set-TextValue $ws.Range("D22") "5.252"
$ws.Range("E22").Value = "  -2.61%  "
Set-TextValue $ws.Range("D23") "10.95"
$ws.Range("E23").Value = "  -2.00%  "
Set-TextValue $ws.Range("D24") "2.105"
$ws.Range("E24").Value = "  +1.73%  "
Set-TextValue $ws.Range("D25") "156.62"
$ws.Range("E25").Value = "  +0.65%  "
Set-TextValue $ws.Range("D26") "19.36"
$ws.Range("E26").Value = "  -2.13%  "
Set-TextValue $ws.Range("D27") "5.984"
$ws.Range("E27").Value = "  +1.26%  "
Set-TextValue $ws.Range("D28") "117.56"
$ws.Range("E28").Value = "  -2.87%  "
Set-TextValue $ws.Range("D29") "1.906"
$ws.Range("E29").Value = "  -5.99%  "
Set-TextValue $ws.Range("D30") "0.09344"
$ws.Range("E30").Value = "  -0.57%  "
Set-TextValue $ws.Range("D31") "0.9008"
$ws.Range("E31").Value = "  -4.54%  "
Set-TextValue $ws.Range("D32") "5.259"
$ws.Range("E32").Value = "  -2.32%  "
Set-TextValue $ws.Range("D33") "1.318"
$ws.Range("E33").Value = "  -3.47%  "
Set-TextValue $ws.Range("D34") "3.252"
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.05737"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D36") "1.168"
$ws.Range("E36").Value = "  -1.10%  "
Set-TextValue $ws.Range("D37") "0.02070"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -0.20%  "
Set-TextValue $ws.Range("D39") "7.641"
$ws.Range("E39").Value = "  -6.58%  "
Set-TextValue $ws.Range("D40") "0.5650"
$ws.Range("E40").Value = "  -3.66%  "
Set-TextValue $ws.Range("D41") "0.1764"
$ws.Range("E41").Value = "  -3.22%  "
Set-TextValue $ws.Range("D42") "9.714"
$ws.Range("E42").Value = "  -4.08%  "
Set-TextValue $ws.Range("D43") "2.236"
$ws.Range("E43").Value = "  -3.96%  "
Set-TextValue $ws.Range("D44") "11.91"
$ws.Range("E44").Value = "  -0.56%  "
Set-TextValue $ws.Range("D45") "0.5339"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("E47").Value = "  -2.72%  "
Set-TextValue $ws.Range("D48") "2.545"
$ws.Range("E48").Value = "  +1.55%  "
Set-TextValue $ws.Range("D49") "112.59"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("E50").Value = "  -6.33%  "
Set-TextValue $ws.Range("D51") "70.69"
$ws.Range("E51").Value = "  -1.36%  "
